$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$ws.Range("B43").Value = "p<=15`np is definied in executive decree 37121"
$ws.Range("B44").Value = "15>p=<30`np is definied in executive decree 37121"
$ws.Range("B45").Value = "30<p<=100`np is definied in executive decree 37121"
$ws.Range("B46").Value = "p > 100`np is definied in executive decree 37121"
